# SOC5060_Framework_Monti_MakingCommunity_2019v00.pptx edit
#
# 1) The "datetimeFigureOut" date placeholder (auto date field) on the
#    slide master and on every slide layout gets its cached display
#    text bumped from 10/9/2019 to 10/17/2019.
# 2) On slide 1 (the "Government/Ethnic/Business/Consumer Approach"
#    diagram), all of the loose shapes that make up the diagram get
#    collected into a single new group named "Group 1".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the "Update automatically" date placeholder everywhere it
#    appears (slide master + all slide layouts) so it reads 10/17/2019.
# ---------------------------------------------------------------------

function Update-DatePlaceholders {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                # ppPlaceholderDate
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder) {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes "10/17/2019"

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholders $layout.Shapes "10/17/2019"
}

# ---------------------------------------------------------------------
# 2) Group every shape on slide 1 into a single new "Group 1" group.
# ---------------------------------------------------------------------

$slide1 = $p.Slides.Item(1)
$shapeCount = $slide1.Shapes.Count

$allIndexes = 1..$shapeCount
$range = $slide1.Shapes.Range($allIndexes)
$group = $range.Group()
$group.Name = "Group 1"
